# Update "想去人数" (interest count) values in F column on sheet "展览" (index 1)
# and sheet "全部类型" (index 4), matching the gh-pages data refresh commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item(1)   # 展览
$wsExhibition.Cells.Item(2, 6).Value = 15235   # was 15234
$wsExhibition.Cells.Item(5, 6).Value = 1594   # was 1590
$wsExhibition.Cells.Item(8, 6).Value = 142   # was 141
$wsExhibition.Cells.Item(9, 6).Value = 258   # was 256
$wsExhibition.Cells.Item(10, 6).Value = 8460   # was 8439
$wsExhibition.Cells.Item(11, 6).Value = 997   # was 996
$wsExhibition.Cells.Item(13, 6).Value = 17   # was 16
$wsExhibition.Cells.Item(16, 6).Value = 76   # was 70
$wsExhibition.Cells.Item(19, 6).Value = 9085   # was 9062
$wsExhibition.Cells.Item(20, 6).Value = 148   # was 147
$wsExhibition.Cells.Item(21, 6).Value = 90   # was 89
$wsExhibition.Cells.Item(22, 6).Value = 205   # was 204
$wsExhibition.Cells.Item(23, 6).Value = 164   # was 163
$wsExhibition.Cells.Item(24, 6).Value = 336   # was 331
$wsExhibition.Cells.Item(25, 6).Value = 5930   # was 5888
$wsExhibition.Cells.Item(26, 6).Value = 1033   # was 1032
$wsExhibition.Cells.Item(27, 6).Value = 46   # was 44
$wsExhibition.Cells.Item(29, 6).Value = 88   # was 86

$wsAllTypes = $wb.Worksheets.Item(4)   # 全部类型
$wsAllTypes.Cells.Item(2, 6).Value = 15235   # was 15234
$wsAllTypes.Cells.Item(5, 6).Value = 1594   # was 1590
$wsAllTypes.Cells.Item(9, 6).Value = 142   # was 141
$wsAllTypes.Cells.Item(10, 6).Value = 258   # was 256
$wsAllTypes.Cells.Item(11, 6).Value = 8460   # was 8439
$wsAllTypes.Cells.Item(12, 6).Value = 997   # was 996
$wsAllTypes.Cells.Item(14, 6).Value = 17   # was 16
$wsAllTypes.Cells.Item(17, 6).Value = 76   # was 70
$wsAllTypes.Cells.Item(22, 6).Value = 9086   # was 9062
$wsAllTypes.Cells.Item(23, 6).Value = 148   # was 147
$wsAllTypes.Cells.Item(24, 6).Value = 90   # was 89
$wsAllTypes.Cells.Item(25, 6).Value = 205   # was 204
$wsAllTypes.Cells.Item(26, 6).Value = 164   # was 163
$wsAllTypes.Cells.Item(27, 6).Value = 336   # was 331
$wsAllTypes.Cells.Item(28, 6).Value = 5930   # was 5888
$wsAllTypes.Cells.Item(29, 6).Value = 1033   # was 1032
$wsAllTypes.Cells.Item(30, 6).Value = 46   # was 44
$wsAllTypes.Cells.Item(32, 6).Value = 88   # was 86
